# HYW Item Balance.xlsx edit script
# - Renames "Pikes & Halberds & Voulges" sheet to "Pikes & Halberds"
# - Fills in item stat data (and derived formula columns) on the
#   "Swords & Daggers" and "Pikes & Halberds" sheets
# - Updates the active sheet / selection state to match the authored edit

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet rename: "Pikes & Halberds & Voulges" -> "Pikes & Halberds"
# ---------------------------------------------------------------------
$wsPikes = $wb.Worksheets.Item(5)
$wsPikes.Name = "Pikes & Halberds"

$wsSwords = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet3 "Swords & Daggers": a few existing F-column values tweaked
# ---------------------------------------------------------------------
$wsSwords.Range("F2").Value = 105
$wsSwords.Range("F3").Value = 108
$wsSwords.Range("F4").Value = 105

# ---------------------------------------------------------------------
# Sheet3 "Swords & Daggers": fill in rows 6-16 (J = SUM(F:I), not shared)
# ---------------------------------------------------------------------
$rowsSum = @(
    @{R=6;  C=325; D=1.2; E=0; F=103; H=34; I=19},
    @{R=7;  C=488; D=1.4; E=0; F=100; H=28; I=26},
    @{R=8;  C=572; D=1.5; E=0; F=99;  H=29; I=25},
    @{R=9;  C=488; D=1.4; E=0; F=100; H=29; I=25},
    @{R=10; C=572; D=1.5; E=0; F=99;  H=30; I=24},
    @{R=11; C=400; D=1.3; E=0; F=102; H=28; I=26},
    @{R=12; C=488; D=1.4; E=0; F=100; H=29; I=25},
    @{R=13; C=572; D=1.5; E=0; F=99;  H=30; I=24},
    @{R=14; C=506; D=1.4; E=0; F=100; H=32; I=28},
    @{R=15; C=402; D=1.3; E=0; F=103; H=28; I=33},
    @{R=16; C=406; D=1.3; E=0; F=102; H=26; I=31}
)

foreach ($row in $rowsSum) {
    $r = $row.R
    $wsSwords.Cells.Item($r, 3).Value = $row.C
    $wsSwords.Cells.Item($r, 4).Value = $row.D
    $wsSwords.Cells.Item($r, 5).Value = $row.E
    $wsSwords.Cells.Item($r, 6).Value = $row.F
    $wsSwords.Cells.Item($r, 8).Value = $row.H
    $wsSwords.Cells.Item($r, 9).Value = $row.I
    $wsSwords.Cells.Item($r, 10).Formula = "=SUM(F$r`:I$r)"
    $wsSwords.Cells.Item($r, 11).Formula = "=J$r*D$r^2"
}

# ---------------------------------------------------------------------
# Sheet3 "Swords & Daggers": fill in rows 18-34 (J = SUM(D:I), shared)
# ---------------------------------------------------------------------
$rowsDI = @(
    @{R=18; C=415; D=1.3; E=0;    F=102; H=36; I=21},
    @{R=19; C=487; D=1.4; E=0;    F=101; H=30; I=23},
    @{R=20; C=582; D=1.5; E=0;    F=99;  H=30; I=25},
    @{R=21; C=508; D=1.4; E=0;    F=100; H=27; I=34},
    @{R=22; C=514; D=1.4; E=0;    F=100; H=28; I=34},
    @{R=23; C=582; D=1.5; E=0;    F=99;  H=28; I=30},
    @{R=24; C=577; D=1.5; E=0;    F=100; H=29; I=27},
    @{R=25; C=670; D=1.6; E=0;    F=97;  H=28; I=28},
    @{R=26; C=491; D=1.4; E=0;    F=101; H=29; I=26},
    @{R=27; C=491; D=1.4; E=$null;F=101; H=30; I=25},
    @{R=28; C=406; D=1.3; E=0;    F=102; H=30; I=24},
    @{R=29; C=497; D=1.4; E=0;    F=101; H=27; I=32},
    @{R=30; C=505; D=1.4; E=0;    F=100; H=31; I=29},
    @{R=31; C=570; D=1.5; E=0;    F=99;  H=31; I=20},
    @{R=32; C=406; D=1.3; E=0;    F=102; H=32; I=22},
    @{R=33; C=499; D=1.4; E=0;    F=100; H=29; I=27},
    @{R=34; C=577; D=1.5; E=0;    F=99;  H=30; I=24}
)

foreach ($row in $rowsDI) {
    $r = $row.R
    $wsSwords.Cells.Item($r, 3).Value = $row.C
    $wsSwords.Cells.Item($r, 4).Value = $row.D
    if ($null -ne $row.E) {
        $wsSwords.Cells.Item($r, 5).Value = $row.E
    }
    $wsSwords.Cells.Item($r, 6).Value = $row.F
    $wsSwords.Cells.Item($r, 8).Value = $row.H
    $wsSwords.Cells.Item($r, 9).Value = $row.I
    $wsSwords.Cells.Item($r, 11).Formula = "=J$r*D$r^2"
}

# J18:J34 is one shared formula covering the whole block
$wsSwords.Range("J18:J34").Formula = "=SUM(D18:I18)"

# ---------------------------------------------------------------------
# Sheet3 view state: no longer the selected tab, selection moved
# ---------------------------------------------------------------------
$wsSwords.Range("F4").Select()

# ---------------------------------------------------------------------
# Sheet5 "Pikes & Halberds": fill in rows 2-10 (J = SUM(D:I)/1.2)
# ---------------------------------------------------------------------
$rowsDiv = @(
    @{R=2;  C=278; D=2;   E=0;  F=85; G=185; H=27; I=34},
    @{R=3;  C=283; D=2;   E=0;  F=85; G=188; H=29; I=35},
    @{R=4;  C=280; D=2;   E=0;  F=85; G=188; H=28; I=33},
    @{R=5;  C=275; D=2;   E=0;  F=85; G=187; H=22; I=36},
    @{R=6;  C=273; D=2;   E=0;  F=85; G=181; H=24; I=35},
    @{R=7;  C=275; D=2;   E=0;  F=85; G=193; H=19; I=31},
    @{R=8;  C=492; D=3.5; E=11; F=78; G=450; H=18; I=30},
    @{R=9;  C=341; D=3;   E=9;  F=82; G=255; H=27; I=33},
    @{R=10; C=333; D=3;   E=9;  F=82; G=246; H=24; I=35}
)

foreach ($row in $rowsDiv) {
    $r = $row.R
    $wsPikes.Cells.Item($r, 3).Value = $row.C
    $wsPikes.Cells.Item($r, 4).Value = $row.D
    $wsPikes.Cells.Item($r, 5).Value = $row.E
    $wsPikes.Cells.Item($r, 6).Value = $row.F
    $wsPikes.Cells.Item($r, 7).Value = $row.G
    $wsPikes.Cells.Item($r, 8).Value = $row.H
    $wsPikes.Cells.Item($r, 9).Value = $row.I
    $wsPikes.Cells.Item($r, 10).Formula = "=SUM(D$r`:I$r)/1.2"
}

# ---------------------------------------------------------------------
# Sheet5 "Pikes & Halberds": fill in rows 12-25 (J = SUM(D:I)*1.2, except
# row 16 which uses *1.1)
# ---------------------------------------------------------------------
$rowsMul = @(
    @{R=12; C=406; D=3.5; E=9; F=86; G=170; H=39; I=31; Mult=1.2},
    @{R=13; C=428; D=3.5; E=9; F=82; G=193; H=37; I=32; Mult=1.2},
    @{R=14; C=423; D=3.6; E=9; F=84; G=180; H=42; I=34; Mult=1.2},
    @{R=15; C=412; D=3.5; E=9; F=84; G=180; H=33; I=34; Mult=1.2},
    @{R=16; C=396; D=3.6; E=9; F=82; G=190; H=39; I=36; Mult=1.1},
    @{R=17; C=415; D=3.8; E=9; F=85; G=173; H=41; I=34; Mult=1.2},
    @{R=18; C=428; D=3.7; E=9; F=82; G=188; H=40; I=34; Mult=1.2},
    @{R=19; C=417; D=3.7; E=9; F=85; G=174; H=42; I=34; Mult=1.2},
    @{R=20; C=435; D=3.6; E=9; F=82; G=190; H=43; I=35; Mult=1.2},
    @{R=21; C=436; D=3.6; E=9; F=82; G=193; H=41; I=35; Mult=1.2},
    @{R=22; C=458; D=3.7; E=9; F=81; G=209; H=44; I=35; Mult=1.2},
    @{R=23; C=406; D=3.6; E=9; F=84; G=172; H=43; I=27; Mult=1.2},
    @{R=24; C=441; D=3.9; E=9; F=82; G=192; H=45; I=36; Mult=1.2},
    @{R=25; C=419; D=4;   E=9; F=82; G=175; H=46; I=33; Mult=1.2}
)

foreach ($row in $rowsMul) {
    $r = $row.R
    $wsPikes.Cells.Item($r, 3).Value = $row.C
    $wsPikes.Cells.Item($r, 4).Value = $row.D
    $wsPikes.Cells.Item($r, 5).Value = $row.E
    $wsPikes.Cells.Item($r, 6).Value = $row.F
    $wsPikes.Cells.Item($r, 7).Value = $row.G
    $wsPikes.Cells.Item($r, 8).Value = $row.H
    $wsPikes.Cells.Item($r, 9).Value = $row.I
    $wsPikes.Cells.Item($r, 10).Formula = "=SUM(D$r`:I$r)*$($row.Mult)"
}

# ---------------------------------------------------------------------
# Sheet5 view state: becomes the selected tab, selection moved
# ---------------------------------------------------------------------
$wsPikes.Activate()
$wsPikes.Range("D32").Select()
